# Replace English control-file text/sheet names with the Danish versions.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the three worksheet tabs (English -> Danish) -----------------
$sheetNameMap = @{
    "19 - 03. Control floating laye" = "19 - 03. Kontrol flydelag"
    "23 - 04. Feeding documentation" = "23 - 04. Foderindlægssedler"
    "77 - 20. Task completed"        = "77 - 20. Arbejdsopgave udført"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($sheetNameMap.ContainsKey($oldName)) {
        $ws.Name = $sheetNameMap[$oldName]
    }
}

# --- 2. Translate cell text content (headers + data) from English to Danish -
$textMap = @{
    "Property"                                  = "Ejendom"
    "Created At"                                 = "Dato"
    "Done By"                                    = "Udført af"
    "Item name"                                  = "Område"
    "Floating layer OK"                          = "Flydelag OK"
    "Select reason for lack of floating layer"   = "Vælg årsag til manglende flydelag"
    "Comment"                                    = "Kommentar"
    "Tjørntved"                                  = "Farm 1"
    "G1: Floating layer"                         = "G1: Flydelag"
    "G2: Floating layer"                         = "G2: Flydelag"
    "Slurry tank empty"                          = "Beholder tom"
    "G3: Floating layer"                         = "G3: Flydelag"
    "Task completed"                             = "Opgave udført"
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value2
            if ($null -ne $val -and $textMap.ContainsKey($val)) {
                $cell.Value = $textMap[$val]
            }
        }
    }
}
